$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 120
$ws.Range("E3").Value = 120
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 110
